# Auto-generated Excel COM-interop script
# Applies updated cryptocurrency price/volume data to sheet1 (cryptos list)
# Commit: "Updated cryptos list on Wed Sep  4 09:29:12 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some 'Price' column values are plain numeric-looking text (e.g. '0.999', '21.70')
# that must stay literal text (matching the source feed's formatting, incl. trailing
# zeros). Force those specific cells to Text format first so assigning the string
# value does not get auto-converted to a number by Excel's smart input parsing.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D10', 'D13', 'D16', 'D19', 'D20', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D35', 'D36', 'D37', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D50')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '56.557.46'
$ws.Range('E2').Value = '  -3.73%  '
$ws.Range('D3').Value = '2.399.76'
$ws.Range('E3').Value = '  -3.67%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '506.01'
$ws.Range('E5').Value = '  -5.00%  '
$ws.Range('D6').Value = '129.97'
$ws.Range('E6').Value = '  -2.79%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').Value = '0.554'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('D9').Value = '2.398.04'
$ws.Range('E9').Value = '  -4.11%  '
$ws.Range('D10').Value = '0.0966'
$ws.Range('E10').Value = '  -2.58%  '
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('D13').Value = '4.66'
$ws.Range('E13').Value = '  -10.44%  '
$ws.Range('D14').Value = '2.821.93'
$ws.Range('E14').Value = '  -3.82%  '
$ws.Range('D15').Value = '56.526.92'
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').Value = '21.70'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('E17').Value = '  -2.70%  '
$ws.Range('D18').Value = '2.384.86'
$ws.Range('E18').Value = '  -4.14%  '
$ws.Range('D19').Value = '10.24'
$ws.Range('E19').Value = '  -3.10%  '
$ws.Range('D20').Value = '313.27'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D22').Value = '6.25'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '65.70'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('D26').Value = '2.491.69'
$ws.Range('E26').Value = '  -4.39%  '
$ws.Range('D27').Value = '0.380'
$ws.Range('E27').Value = '  -6.56%  '
$ws.Range('D28').Value = '0.151'
$ws.Range('E28').Value = '  -5.06%  '
$ws.Range('D29').Value = '7.27'
$ws.Range('E29').Value = '  -2.24%  '
$ws.Range('D30').Value = '174.49'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.67'
$ws.Range('E31').Value = '  -3.33%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0716'
$ws.Range('E32').Value = '  -5.15%  '
$ws.Range('D33').Value = '6.15'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('E34').Value = '  -5.92%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '17.80'
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('D40').Value = '35.86'
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('D41').Value = '1.44'
$ws.Range('E41').Value = '  -4.30%  '
$ws.Range('D42').Value = '0.786'
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('D43').Value = '132.85'
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('D45').Value = '4.90'
$ws.Range('E45').Value = '  -4.71%  '
$ws.Range('D46').Value = '255.97'
$ws.Range('E46').Value = '  -6.54%  '
$ws.Range('E47').Value = '  -3.31%  '
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').Value = '16.83'
$ws.Range('E50').Value = '  -4.18%  '
$ws.Range('E51').Value = '  -4.44%  '
